$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This template is used to generate "Teilnehmerlisten" (participant lists) for
# a "Freizeit" (camp/event). Replace the static placeholder labels with the
# actual merge-field tokens used by the TN-list generator.

# Header block: event name (was the hard-coded "Freizeitname" label).
$ws.Range("D2").Value = '${bezeichnung}'

# Leader line: event leader placeholder (was the hard-coded "Freizeitleiter").
$ws.Range("D3").Value = '????'

# Date range: was literal sample dates, now merge-field tokens.
$ws.Range("F3").Value = '${begin.german}'
$ws.Range("G3").Value = '${ende.german}'

# Reflect the editor's last selection being near the edited cells.
$ws.Range("D4").Select()
